# Add data for 2022-05-04 update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet and update the "through" date label
$ws.Name = "Through 2022-04-26"
$ws.Range("I1").Value = "2022 (through 04-26)"

# Update data values per the diff
$ws.Range("I5").Value = 108
$ws.Range("H13").Value = 205
$ws.Range("H14").Value = 1852
$ws.Range("I14").Value = 543
